$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: replace full distributor part-numbers with plain item numbers ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Re-home A3/A4 onto the same plain font used by the "Manufactur Part#" column
# so that the font slot they used to own is freed up for the new hyperlink style.
$ws.Range("A3:A4").Font.Color = $ws.Range("C2").Font.Color
$ws.Range("A3:A4").Font.Name = $ws.Range("C2").Font.Name

# --- Column F: currency changed from USD to CAD ---
$ws.Range("F2").Value = "CAD"
$ws.Range("F3").Value = "CAD"
$ws.Range("F4").Value = "CAD"

# --- Column G: totals stored as plain numeric values (no more "x * 1.3" text) ---
$ws.Range("G2").Value = 1.35
$ws.Range("G3").Value = 2.32
$ws.Range("G4").Value = 4.79

# --- Column C: turn the manufacturer part numbers into hyperlinks ---
$ws.Range("C2:C4").ClearFormats()
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.mouser.com/ProductDetail/538-51216-0200")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.mouser.com/ProductDetail/538-51217-0205")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.mouser.com/ProductDetail/538-55755-0219")

# --- Narrower item-number column now that it only holds small integers ---
$ws.Columns.Item(1).ColumnWidth = 6.166666666666667

# --- Misc view / print settings ---
$ws.PageSetup.Orientation = 1
[void]$ws.Range("C4").Select()
